$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add note to J6 explaining that tax and fort are halved from total
$ws.Range("J6").Value = "Tax and fort are halved from total."

# New "Lorathi City" block (rows 35-45), following the same pattern as the
# other culture blocks above it.
$ws.Range("A35").Value = "Lorathi City"
$ws.Range("B35").Value = 0.45
$ws.Range("C35").Value = -0.2
$ws.Range("D35").Value = 10
$ws.Range("E35").Value = 15
$ws.Range("I35").Value = 75

$ws.Range("A36").Value = 2
$ws.Range("B36").Value = -0.43
$ws.Range("C36").Value = 0.7

$ws.Range("A37").Value = 3
$ws.Range("B37").Value = -0.43
$ws.Range("C37").Value = -0.9

$ws.Range("A38").Value = 4
$ws.Range("B38").Value = 0.18
$ws.Range("C38").Value = 0.5

$ws.Range("A39").Value = 5
$ws.Range("B39").Value = 0.09
$ws.Range("C39").Value = 0.2

$ws.Range("A40").Value = 6
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = -0.7

$ws.Range("A41").Value = 7
$ws.Range("B41").Value = 0.29
$ws.Range("C41").Value = -0.9

$ws.Range("A42").Value = 8
$ws.Range("B42").Value = -0.03
$ws.Range("C42").Value = -0.5

$ws.Range("A43").Value = 9
$ws.Range("B43").Value = -0.37
$ws.Range("C43").Value = 0.6

$ws.Range("A44").Value = 10
$ws.Range("B44").Value = -0.08
$ws.Range("C44").Value = 0.5

$ws.Range("A45").Value = "Total"
$ws.Range("B45").Formula = "=SUM(B35:B44)"
$ws.Range("C45").Formula = "=SUM(C35:C44)"

# Update the view to match the committed state
$ws.Range("E36").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
